{"js": "// Replace the 100 arithmetic-expression cells of the single table in the\n// document with their updated values, in document order (row-major,\n// left-to-right). A positional mapping is required because a handful of\n// the old expressions repeat verbatim (e.g. \"58+12=70\" appears twice) but\n// map to two different new expressions depending on where they sit in the\n// table, so a global text find/replace would be ambiguous.\nconst replacements = [\n  [\"94-90=4\", \"69-45=24\"], [\"63+13=76\", \"46+10=56\"], [\"30-0=30\", \"74-10=64\"],\n  [\"63-5=58\", \"41-30=11\"], [\"57+23=80\", \"32-24=8\"], [\"20+3=23\", \"86-50=36\"],\n  [\"53+33=86\", \"54-48=6\"], [\"16+14=30\", \"79-35=44\"], [\"37+37=74\", \"81-36=45\"],\n  [\"40+35=75\", \"90-68=22\"], [\"38-31=7\", \"12+10=22\"], [\"16+36=52\", \"34-16=18\"],\n  [\"8+37=45\", \"87-22=65\"], [\"28+49=77\", \"69-34=35\"], [\"5+38=43\", \"36+44=80\"],\n  [\"14+5=19\", \"1+48=49\"], [\"21+46=67\", \"96-61=35\"], [\"8-4=4\", \"29+54=83\"],\n  [\"58-29=29\", \"13+46=59\"], [\"73-42=31\", \"74+24=98\"], [\"23+0=23\", \"23+65=88\"],\n  [\"67+17=84\", \"71+4=75\"], [\"67-12=55\", \"52+15=67\"], [\"16+11=27\", \"49-40=9\"],\n  [\"58+12=70\", \"96-50=46\"], [\"13+5=18\", \"14+19=33\"], [\"0+80=80\", \"54-16=38\"],\n  [\"13-3=10\", \"56-21=35\"], [\"80+9=89\", \"55-23=32\"], [\"63-63=0\", \"17-2=15\"],\n  [\"37+51=88\", \"52+23=75\"], [\"11+19=30\", \"88+3=91\"], [\"98-59=39\", \"22+43=65\"],\n  [\"64+3=67\", \"74+23=97\"], [\"21+1=22\", \"70-10=60\"], [\"26+35=61\", \"41+7=48\"],\n  [\"53+31=84\", \"60+6=66\"], [\"97-11=86\", \"40+18=58\"], [\"30+5=35\", \"34+57=91\"],\n  [\"46+35=81\", \"8+28=36\"], [\"91-42=49\", \"46+42=88\"], [\"65-32=33\", \"88-49=39\"],\n  [\"18-15=3\", \"86+8=94\"], [\"11+16=27\", \"10-9=1\"], [\"40-21=19\", \"19-7=12\"],\n  [\"27+15=42\", \"48+49=97\"], [\"97-34=63\", \"44+42=86\"], [\"55-15=40\", \"31+50=81\"],\n  [\"68-24=44\", \"21+18=39\"], [\"48+22=70\", \"22+17=39\"], [\"87+4=91\", \"72+16=88\"],\n  [\"55+2=57\", \"33-13=20\"], [\"49-46=3\", \"51-14=37\"], [\"7+55=62\", \"38+52=90\"],\n  [\"28-3=25\", \"52+4=56\"], [\"45-17=28\", \"72-37=35\"], [\"48-33=15\", \"50-48=2\"],\n  [\"78-59=19\", \"44+38=82\"], [\"51+0=51\", \"82-74=8\"], [\"63+20=83\", \"34+26=60\"],\n  [\"10+44=54\", \"1+70=71\"], [\"11+86=97\", \"70+1=71\"], [\"21-9=12\", \"98-80=18\"],\n  [\"0+13=13\", \"42-32=10\"], [\"76+3=79\", \"76-55=21\"], [\"10+82=92\", \"49+8=57\"],\n  [\"49+12=61\", \"57-41=16\"], [\"27-22=5\", \"20+71=91\"], [\"92-7=85\", \"32+7=39\"],\n  [\"86-21=65\", \"3+54=57\"], [\"8+34=42\", \"96-15=81\"], [\"16+44=60\", \"67+0=67\"],\n  [\"95-66=29\", \"6+34=40\"], [\"75-17=58\", \"9+46=55\"], [\"71-20=51\", \"37-29=8\"],\n  [\"46-39=7\", \"96-13=83\"], [\"45+21=66\", \"77-4=73\"], [\"43-38=5\", \"55-9=46\"],\n  [\"46+23=69\", \"29+5=34\"], [\"87-11=76\", \"56-14=42\"], [\"18+78=96\", \"88-19=69\"],\n  [\"94-15=79\", \"44-26=18\"], [\"36+32=68\", \"7+46=53\"], [\"89-40=49\", \"18+6=24\"],\n  [\"76-37=39\", \"85-29=56\"], [\"63-40=23\", \"46+18=64\"], [\"29-26=3\", \"8-1=7\"],\n  [\"85-40=45\", \"29+58=87\"], [\"51-37=14\", \"33+58=91\"], [\"45-12=33\", \"6+84=90\"],\n  [\"18+53=71\", \"50-13=37\"], [\"52+27=79\", \"61+0=61\"], [\"34+55=89\", \"11+22=33\"],\n  [\"28+44=72\", \"55+9=64\"], [\"73-63=10\", \"87-65=22\"], [\"95-21=74\", \"67-40=27\"],\n  [\"73-1=72\", \"88-27=61\"], [\"58+12=70\", \"37+42=79\"], [\"90-22=68\", \"39+55=94\"],\n  [\"70-11=59\", \"22+63=85\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document but found none.\");\n}\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load every row's cells up front.\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\n// Flatten all cells in row-major (document) order.\nconst cells = [];\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    cells.push(cell);\n  }\n}\n\nif (cells.length !== replacements.length) {\n  throw new Error(\n    \"Expected \" + replacements.length + \" table cells but found \" + cells.length\n  );\n}\n\n// Load current text of every cell so we can sanity-check before writing.\nfor (const cell of cells) {\n  cell.load(\"value\");\n}\nawait context.sync();\n\nfor (let i = 0; i < cells.length; i++) {\n  const [oldText, newText] = replacements[i];\n  const cell = cells[i];\n  if (cell.value !== oldText) {\n    throw new Error(\n      \"Cell \" + i + \" expected '\" + oldText + \"' but found '\" + cell.value + \"'\"\n    );\n  }\n  cell.value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-expression cells of the single table in the\n# document with their updated values, in document order (row-major,\n# left-to-right, matching Cell(row, col) iteration order). A positional\n# mapping is required because a handful of the old expressions repeat\n# verbatim (e.g. \"58+12=70\" appears twice) but map to two different new\n# expressions depending on where they sit in the table, so a global\n# text find/replace would be ambiguous.\n$replacements = @(\n  @(\"94-90=4\", \"69-45=24\"),\n  @(\"63+13=76\", \"46+10=56\"),\n  @(\"30-0=30\", \"74-10=64\"),\n  @(\"63-5=58\", \"41-30=11\"),\n  @(\"57+23=80\", \"32-24=8\"),\n  @(\"20+3=23\", \"86-50=36\"),\n  @(\"53+33=86\", \"54-48=6\"),\n  @(\"16+14=30\", \"79-35=44\"),\n  @(\"37+37=74\", \"81-36=45\"),\n  @(\"40+35=75\", \"90-68=22\"),\n  @(\"38-31=7\", \"12+10=22\"),\n  @(\"16+36=52\", \"34-16=18\"),\n  @(\"8+37=45\", \"87-22=65\"),\n  @(\"28+49=77\", \"69-34=35\"),\n  @(\"5+38=43\", \"36+44=80\"),\n  @(\"14+5=19\", \"1+48=49\"),\n  @(\"21+46=67\", \"96-61=35\"),\n  @(\"8-4=4\", \"29+54=83\"),\n  @(\"58-29=29\", \"13+46=59\"),\n  @(\"73-42=31\", \"74+24=98\"),\n  @(\"23+0=23\", \"23+65=88\"),\n  @(\"67+17=84\", \"71+4=75\"),\n  @(\"67-12=55\", \"52+15=67\"),\n  @(\"16+11=27\", \"49-40=9\"),\n  @(\"58+12=70\", \"96-50=46\"),\n  @(\"13+5=18\", \"14+19=33\"),\n  @(\"0+80=80\", \"54-16=38\"),\n  @(\"13-3=10\", \"56-21=35\"),\n  @(\"80+9=89\", \"55-23=32\"),\n  @(\"63-63=0\", \"17-2=15\"),\n  @(\"37+51=88\", \"52+23=75\"),\n  @(\"11+19=30\", \"88+3=91\"),\n  @(\"98-59=39\", \"22+43=65\"),\n  @(\"64+3=67\", \"74+23=97\"),\n  @(\"21+1=22\", \"70-10=60\"),\n  @(\"26+35=61\", \"41+7=48\"),\n  @(\"53+31=84\", \"60+6=66\"),\n  @(\"97-11=86\", \"40+18=58\"),\n  @(\"30+5=35\", \"34+57=91\"),\n  @(\"46+35=81\", \"8+28=36\"),\n  @(\"91-42=49\", \"46+42=88\"),\n  @(\"65-32=33\", \"88-49=39\"),\n  @(\"18-15=3\", \"86+8=94\"),\n  @(\"11+16=27\", \"10-9=1\"),\n  @(\"40-21=19\", \"19-7=12\"),\n  @(\"27+15=42\", \"48+49=97\"),\n  @(\"97-34=63\", \"44+42=86\"),\n  @(\"55-15=40\", \"31+50=81\"),\n  @(\"68-24=44\", \"21+18=39\"),\n  @(\"48+22=70\", \"22+17=39\"),\n  @(\"87+4=91\", \"72+16=88\"),\n  @(\"55+2=57\", \"33-13=20\"),\n  @(\"49-46=3\", \"51-14=37\"),\n  @(\"7+55=62\", \"38+52=90\"),\n  @(\"28-3=25\", \"52+4=56\"),\n  @(\"45-17=28\", \"72-37=35\"),\n  @(\"48-33=15\", \"50-48=2\"),\n  @(\"78-59=19\", \"44+38=82\"),\n  @(\"51+0=51\", \"82-74=8\"),\n  @(\"63+20=83\", \"34+26=60\"),\n  @(\"10+44=54\", \"1+70=71\"),\n  @(\"11+86=97\", \"70+1=71\"),\n  @(\"21-9=12\", \"98-80=18\"),\n  @(\"0+13=13\", \"42-32=10\"),\n  @(\"76+3=79\", \"76-55=21\"),\n  @(\"10+82=92\", \"49+8=57\"),\n  @(\"49+12=61\", \"57-41=16\"),\n  @(\"27-22=5\", \"20+71=91\"),\n  @(\"92-7=85\", \"32+7=39\"),\n  @(\"86-21=65\", \"3+54=57\"),\n  @(\"8+34=42\", \"96-15=81\"),\n  @(\"16+44=60\", \"67+0=67\"),\n  @(\"95-66=29\", \"6+34=40\"),\n  @(\"75-17=58\", \"9+46=55\"),\n  @(\"71-20=51\", \"37-29=8\"),\n  @(\"46-39=7\", \"96-13=83\"),\n  @(\"45+21=66\", \"77-4=73\"),\n  @(\"43-38=5\", \"55-9=46\"),\n  @(\"46+23=69\", \"29+5=34\"),\n  @(\"87-11=76\", \"56-14=42\"),\n  @(\"18+78=96\", \"88-19=69\"),\n  @(\"94-15=79\", \"44-26=18\"),\n  @(\"36+32=68\", \"7+46=53\"),\n  @(\"89-40=49\", \"18+6=24\"),\n  @(\"76-37=39\", \"85-29=56\"),\n  @(\"63-40=23\", \"46+18=64\"),\n  @(\"29-26=3\", \"8-1=7\"),\n  @(\"85-40=45\", \"29+58=87\"),\n  @(\"51-37=14\", \"33+58=91\"),\n  @(\"45-12=33\", \"6+84=90\"),\n  @(\"18+53=71\", \"50-13=37\"),\n  @(\"52+27=79\", \"61+0=61\"),\n  @(\"34+55=89\", \"11+22=33\"),\n  @(\"28+44=72\", \"55+9=64\"),\n  @(\"73-63=10\", \"87-65=22\"),\n  @(\"95-21=74\", \"67-40=27\"),\n  @(\"73-1=72\", \"88-27=61\"),\n  @(\"58+12=70\", \"37+42=79\"),\n  @(\"90-22=68\", \"39+55=94\"),\n  @(\"70-11=59\", \"22+63=85\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\nif (($rowCount * $colCount) -ne $replacements.Count) {\n  throw \"Expected $($replacements.Count) table cells but found $($rowCount * $colCount)\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    $pair = $replacements[$i]\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $cell = $t.Cell($r, $c)\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $oldText) {\n      throw \"Cell ($r,$c) expected '$oldText' but found '$current'\"\n    }\n    $cell.Range.Text = $newText\n\n    $i = $i + 1\n  }\n}\n"}
